$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix comma-separated names to use periods ---
$ws.Range("E45").Value = "FERNANDEZ MARIO H. GALLICET OSCAR M"
$ws.Range("E61").Value = "FERNANDEZ MARIO H. GALLICET OSCAR M"
$ws.Range("E46").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("F46").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("E47").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"

# --- Fix number formatting: remove thousands separator "." and use "." as decimal (was ",") ---
$ws.Range("H2").Value = "'1550.00"
$ws.Range("H3").Value = "'500.00"
$ws.Range("H4").Value = "'24999.00"
$ws.Range("H5").Value = "'189.92"
$ws.Range("H6").Value = "'11280.00"
$ws.Range("H7").Value = "'22169.00"
$ws.Range("H8").Value = "'26824.40"
$ws.Range("H9").Value = "'271.66"
$ws.Range("H10").Value = "'105.73"
$ws.Range("H11").Value = "'2515.50"
$ws.Range("H12").Value = "'5600.57"
$ws.Range("H13").Value = "'333.20"
$ws.Range("H14").Value = "'698.70"
$ws.Range("H15").Value = "'6608.79"
$ws.Range("H16").Value = "'960.00"
$ws.Range("H17").Value = "'80.00"
$ws.Range("H18").Value = "'413.00"
$ws.Range("H19").Value = "'950.00"
$ws.Range("H20").Value = "'137.00"
$ws.Range("H21").Value = "'1172.90"
$ws.Range("H22").Value = "'116.00"
$ws.Range("H23").Value = "'129.10"
$ws.Range("H24").Value = "'9.00"
$ws.Range("H25").Value = "'219.00"
$ws.Range("H26").Value = "'590.00"
$ws.Range("H27").Value = "'2183.37"
$ws.Range("H28").Value = "'630.00"
$ws.Range("H29").Value = "'739.04"
$ws.Range("H30").Value = "'12.00"
$ws.Range("H31").Value = "'264.00"
$ws.Range("H32").Value = "'764.70"
$ws.Range("H33").Value = "'150.00"
$ws.Range("H34").Value = "'6.50"
$ws.Range("H35").Value = "'585.40"
$ws.Range("H36").Value = "'226.38"
$ws.Range("H37").Value = "'510.00"
$ws.Range("H38").Value = "'464.00"
$ws.Range("H39").Value = "'130.00"
$ws.Range("H40").Value = "'1107.15"
$ws.Range("H41").Value = "'4386.00"
$ws.Range("H42").Value = "'2990.00"
$ws.Range("H43").Value = "'10853.20"
$ws.Range("H44").Value = "'203.00"
$ws.Range("H45").Value = "'269.00"
$ws.Range("H46").Value = "'228.60"
$ws.Range("H47").Value = "'550.00"
$ws.Range("H48").Value = "'108.90"
$ws.Range("H49").Value = "'340.00"
$ws.Range("H50").Value = "'499233.27"
$ws.Range("H51").Value = "'110472.98"
$ws.Range("H52").Value = "'64.00"
$ws.Range("H53").Value = "'52.50"
$ws.Range("H54").Value = "'409.19"
$ws.Range("H55").Value = "'3213.00"
$ws.Range("H56").Value = "'72.00"
$ws.Range("H57").Value = "'1809.50"
$ws.Range("H58").Value = "'1690.00"
$ws.Range("H59").Value = "'233.70"
$ws.Range("H60").Value = "'373.00"
$ws.Range("H61").Value = "'2631.60"
$ws.Range("H62").Value = "'386.00"
$ws.Range("H63").Value = "'28.00"
$ws.Range("H64").Value = "'3284.40"
$ws.Range("H65").Value = "'185.50"
$ws.Range("H66").Value = "'104.91"
$ws.Range("H67").Value = "'511.00"
$ws.Range("H68").Value = "'300.00"
$ws.Range("H69").Value = "'47.00"
$ws.Range("H70").Value = "'264.00"
$ws.Range("H71").Value = "'120.00"
$ws.Range("H72").Value = "'757.00"
$ws.Range("H73").Value = "'500.00"
$ws.Range("H74").Value = "'2541.50"
$ws.Range("H75").Value = "'150.00"
$ws.Range("H76").Value = "'18500.00"
$ws.Range("H77").Value = "'8505.87"
$ws.Range("H78").Value = "'200.00"
$ws.Range("H79").Value = "'1250.00"
$ws.Range("H80").Value = "'285.00"
$ws.Range("H81").Value = "'9317.00"
$ws.Range("H82").Value = "'4333.74"
$ws.Range("H83").Value = "'878.00"
$ws.Range("H84").Value = "'4932.72"
$ws.Range("H85").Value = "'600.00"
$ws.Range("H86").Value = "'300.00"
$ws.Range("H87").Value = "'750.00"
$ws.Range("H88").Value = "'120.00"
$ws.Range("H89").Value = "'67.50"
$ws.Range("H90").Value = "'170.00"
$ws.Range("H91").Value = "'1635.00"
$ws.Range("H92").Value = "'1981.00"
$ws.Range("H93").Value = "'174.00"
$ws.Range("H94").Value = "'60.00"
$ws.Range("H95").Value = "'950.00"
$ws.Range("H96").Value = "'450.00"
$ws.Range("H97").Value = "'216.20"
$ws.Range("H98").Value = "'745.00"
$ws.Range("H99").Value = "'258.64"
$ws.Range("H100").Value = "'21636.00"
$ws.Range("H101").Value = "'330.00"
$ws.Range("H102").Value = "'2884.00"
$ws.Range("H103").Value = "'1941.40"
$ws.Range("H104").Value = "'2098.92"
$ws.Range("H105").Value = "'380.00"
$ws.Range("H106").Value = "'19.47"
$ws.Range("H107").Value = "'9100.00"
$ws.Range("H108").Value = "'868.00"
$ws.Range("H109").Value = "'28.00"
$ws.Range("H110").Value = "'723.00"
$ws.Range("H111").Value = "'1762.25"
$ws.Range("H112").Value = "'838.87"
$ws.Range("H113").Value = "'12420.00"
$ws.Range("H114").Value = "'1150.00"
$ws.Range("H115").Value = "'4340.21"
$ws.Range("H116").Value = "'12879.00"
$ws.Range("H117").Value = "'238.00"
$ws.Range("H118").Value = "'847458.12"
$ws.Range("H119").Value = "'1900.00"
$ws.Range("H120").Value = "'1114.26"
$ws.Range("H121").Value = "'4600.00"
